# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values with newly recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    9  = 1
    10 = 5
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 4
    16 = 2
    17 = 1
    18 = 3
    19 = 2
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 4
    34 = 3
    35 = 1
    36 = 1
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 0
    49 = 0
    50 = 2
    51 = 0
    52 = 2
    53 = 2
    54 = 2
    56 = 1
    57 = 3
    58 = 1
    60 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
